$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — copy H1's format (bold,
# centered, bordered header style) onto them, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new columns I (I0) and J (IF), rows 2-37.
$iValues = @(6,7,6,8,7,7,6,7,7,9,5,5,5,8,5,7,6,6,6,7,6,5,7,7,8,9,9,8,7,8,7,8,7,7,2,8)
$jValues = @(6,7,6,8,8,7,6,8,7,9,6,6,6,8,6,7,6,6,6,7,6,6,8,8,8,9,9,8,7,8,7,8,8,8,3,8)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
